$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 183, shifting existing rows
# 183-299 down to 184-300 (dimension grows from A1:T299 to A1:T300).
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row 183 with the new weekly price record
# for "Vega Modelo de Temuco" / Mango (same categorical fields as the
# surrounding rows, new date + price figures).
$ws.Cells.Item(183, 1).Value  = 10
$ws.Cells.Item(183, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(183, 3).Value  = "La Araucanía"
$ws.Cells.Item(183, 4).Value  = 44606
$ws.Cells.Item(183, 5).Value  = 9
$ws.Cells.Item(183, 6).Value  = "Fruta"
$ws.Cells.Item(183, 7).Value  = 100108
$ws.Cells.Item(183, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(183, 9).Value  = 100108002
$ws.Cells.Item(183, 10).Value = "Mango"
$ws.Cells.Item(183, 11).Value = "Sin especificar"
$ws.Cells.Item(183, 12).Value = "Primera"
$ws.Cells.Item(183, 13).Value = 1300
$ws.Cells.Item(183, 14).Value = 7000
$ws.Cells.Item(183, 15).Value = 8000
$ws.Cells.Item(183, 16).Value = 7462
$ws.Cells.Item(183, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(183, 18).Value = "Perú"
$ws.Cells.Item(183, 19).Value = 1866
$ws.Cells.Item(183, 20).Value = 4
